# Milestone 3 Presentation fixes
# 1) Slide 7: finish the truncated sentence about code formalities.
# 2) Theme colors: the deck's active theme (theme2.xml, used by the
#    slide master / presentation) was mistakenly carrying the
#    "Macmorris template" palette; restore the original "Default"
#    color palette on it.

$p = $ppt.ActivePresentation

# --- 1. Fix the truncated sentence on slide 7 -----------------------------
$slide = $p.Slides.Item(7)
$shape = $slide.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

$oldSentence = "Yesterday as a group we all sat down, merged our code, and ensured the code was clean and that the code followed some of the"
$newSentence = "Yesterday as a group we all sat down, merged our code, and ensured the code was clean and that the code followed some of the formalities we went over in class."

$fullText = $tr.Text
$startIdx = $fullText.IndexOf($oldSentence)
if ($startIdx -ge 0) {
    $run = $tr.Characters($startIdx + 1, $oldSentence.Length)
    $run.Text = $newSentence
}

# --- 2. Restore the Default color scheme on the presentation theme -------
function HexToRgbVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + $g * 256 + $b * 65536
}

# Order matches the DrawingML clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$defaultColors = @(
    "000000",
    "FFFFFF",
    "158158",
    "F3F3F3",
    "058DC7",
    "50B432",
    "ED561B",
    "EDEF00",
    "24CBE5",
    "64E572",
    "2200CC",
    "551A8B"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgbVal $defaultColors[$i - 1]
}
